$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal text value without altering its cell style,
# even when the text looks like a number (e.g. "520.25") or like a date/other
# auto-converted type. We temporarily force a Text number format, assign the
# value, then restore the original (Normal) cell style so no stray style index
# is left behind on the cell.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.377.54"
$ws.Range("E2").Value = "  -4.08%  "

Set-TextValue $ws.Range("D3") "2.641.12"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue $ws.Range("D5") "520.25"
$ws.Range("E5").Value = "  -1.21%  "

Set-TextValue $ws.Range("D6") "143.72"
$ws.Range("E6").Value = "  -0.68%  "

Set-TextValue $ws.Range("D8") "0.568"
$ws.Range("E8").Value = "  -1.98%  "

Set-TextValue $ws.Range("D9") "6.67"
$ws.Range("E9").Value = "  -0.28%  "

Set-TextValue $ws.Range("D10") "0.102"
$ws.Range("E10").Value = "  -3.23%  "

Set-TextValue $ws.Range("D11") "0.337"
$ws.Range("E11").Value = "  -1.04%  "

Set-TextValue $ws.Range("D13") "3.105.56"
$ws.Range("E13").Value = "  -2.30%  "

Set-TextValue $ws.Range("D14") "58.373.28"
$ws.Range("E14").Value = "  -4.01%  "

Set-TextValue $ws.Range("D15") "20.76"
$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("E16").Value = "  -1.41%  "

Set-TextValue $ws.Range("D17") "2.647.05"
$ws.Range("E17").Value = "  -7.26%  "

Set-TextValue $ws.Range("D18") "336.74"
$ws.Range("E18").Value = "  -3.23%  "

Set-TextValue $ws.Range("D19") "4.39"
$ws.Range("E19").Value = "  -2.45%  "

Set-TextValue $ws.Range("D20") "10.44"
$ws.Range("E20").Value = "  -1.59%  "

Set-TextValue $ws.Range("D21") "6.28"
$ws.Range("E21").Value = "  -2.08%  "

Set-TextValue $ws.Range("D22") "0.998"
$ws.Range("E22").Value = "  -0.05%  "

Set-TextValue $ws.Range("D23") "64.37"
$ws.Range("E23").Value = "  +0.93%  "

Set-TextValue $ws.Range("D24") "0.423"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("E25").Value = "  -2.07%  "

Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.72%  "

Set-TextValue $ws.Range("D27") "0.0₃0793"
$ws.Range("E27").Value = "  -2.81%  "

Set-TextValue $ws.Range("D28") "7.09"
$ws.Range("E28").Value = "  -3.02%  "

Set-TextValue $ws.Range("D29") "6.60"
$ws.Range("E29").Value = "  -2.48%  "

Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  -1.04%  "

Set-TextValue $ws.Range("D32") "152.56"
$ws.Range("E32").Value = "  +1.39%  "

Set-TextValue $ws.Range("D33") "18.80"
$ws.Range("E33").Value = "  -1.77%  "

Set-TextValue $ws.Range("D34") "4.12"
$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("E35").Value = "  -4.05%  "

Set-TextValue $ws.Range("D36") "0.904"
$ws.Range("E36").Value = "  -4.09%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D37") "36.71"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D38") "0.854"
$ws.Range("E38").Value = "  -2.50%  "

Set-TextValue $ws.Range("D39") "1.44"
$ws.Range("E39").Value = "  -5.03%  "

Set-TextValue $ws.Range("D40") "3.62"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("E41").Value = "  +0.24%  "

Set-TextValue $ws.Range("D42") "0.605"
$ws.Range("E42").Value = "  -1.00%  "

Set-TextValue $ws.Range("D43") "0.0968"
$ws.Range("E43").Value = "  -2.30%  "

Set-TextValue $ws.Range("D44") "268.73"
$ws.Range("E44").Value = "  -5.47%  "

Set-TextValue $ws.Range("D45") "19.33"
$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("E46").Value = "  +1.64%  "

Set-TextValue $ws.Range("D47") "0.0535"
$ws.Range("E47").Value = "  -0.65%  "

Set-TextValue $ws.Range("D48") "2.044.43"
$ws.Range("E48").Value = "  -4.57%  "

Set-TextValue $ws.Range("D49") "4.68"
$ws.Range("E49").Value = "  -2.33%  "

Set-TextValue $ws.Range("D50") "0.0227"
$ws.Range("E50").Value = "  -3.26%  "

Set-TextValue $ws.Range("D51") "18.25"
$ws.Range("E51").Value = "  -5.45%  "
